$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44424
$ws.Range("M2").Value = 25

# Row 3
$ws.Range("D3").Value = 44231
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 25000
$ws.Range("S3").Value = 1250

# Row 4
$ws.Range("D4").Value = 44428
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 24000
$ws.Range("P4").Value = 24000
$ws.Range("S4").Value = 1200

# Row 5
$ws.Range("D5").Value = 44334
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 25000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 25000
$ws.Range("S5").Value = 1250

# Row 6
$ws.Range("D6").Value = 44414
$ws.Range("M6").Value = 15
$ws.Range("N6").Value = 25000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 25000
$ws.Range("S6").Value = 1250

# Row 7
$ws.Range("D7").Value = 44461
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 24000
$ws.Range("P7").Value = 24000
$ws.Range("S7").Value = 1200

# Row 8
$ws.Range("D8").Value = 44466
$ws.Range("M8").Value = 70
$ws.Range("N8").Value = 24000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 24000
$ws.Range("S8").Value = 1200

# Row 9
$ws.Range("D9").Value = 44418

# Row 10
$ws.Range("D10").Value = 44221
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 25000
$ws.Range("O10").Value = 25000
$ws.Range("P10").Value = 25000
$ws.Range("S10").Value = 1250

# Row 11
$ws.Range("D11").Value = 44412
$ws.Range("M11").Value = 20
$ws.Range("N11").Value = 25000
$ws.Range("O11").Value = 25000
$ws.Range("P11").Value = 25000
$ws.Range("S11").Value = 1250

# Row 12
$ws.Range("D12").Value = 44433
$ws.Range("M12").Value = 10
$ws.Range("N12").Value = 24000
$ws.Range("O12").Value = 24000
$ws.Range("P12").Value = 24000
$ws.Range("S12").Value = 1200

# Row 13
$ws.Range("D13").Value = 44392
$ws.Range("M13").Value = 10
$ws.Range("N13").Value = 24000
$ws.Range("O13").Value = 24000
$ws.Range("P13").Value = 24000
$ws.Range("S13").Value = 1200

# Row 14
$ws.Range("D14").Value = 44489

# Row 15
$ws.Range("D15").Value = 44434
$ws.Range("M15").Value = 20

# Row 16
$ws.Range("D16").Value = 44435
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 24000
$ws.Range("O16").Value = 24000
$ws.Range("P16").Value = 24000
$ws.Range("S16").Value = 1200

# Row 17
$ws.Range("D17").Value = 44442
$ws.Range("M17").Value = 25
$ws.Range("N17").Value = 23000
$ws.Range("O17").Value = 23000
$ws.Range("P17").Value = 23000
$ws.Range("S17").Value = 1150

# Row 18
$ws.Range("D18").Value = 44175
$ws.Range("M18").Value = 25
$ws.Range("N18").Value = 23000
$ws.Range("O18").Value = 23000
$ws.Range("P18").Value = 23000
$ws.Range("S18").Value = 1150

# Row 19
$ws.Range("D19").Value = 44475
$ws.Range("M19").Value = 20

# Row 20
$ws.Range("D20").Value = 44235
$ws.Range("M20").Value = 15
$ws.Range("N20").Value = 25000
$ws.Range("O20").Value = 25000
$ws.Range("P20").Value = 25000
$ws.Range("S20").Value = 1250

# Row 21
$ws.Range("D21").Value = 44419
$ws.Range("M21").Value = 40
$ws.Range("N21").Value = 25000
$ws.Range("O21").Value = 25000
$ws.Range("P21").Value = 25000
$ws.Range("S21").Value = 1250

# Row 23
$ws.Range("D23").Value = 44214

# Row 24
$ws.Range("D24").Value = 44390
$ws.Range("M24").Value = 10

# Row 25
$ws.Range("D25").Value = 44420
$ws.Range("M25").Value = 35

# Row 26
$ws.Range("D26").Value = 44398
$ws.Range("M26").Value = 15
$ws.Range("N26").Value = 25000
$ws.Range("O26").Value = 25000
$ws.Range("P26").Value = 25000
$ws.Range("S26").Value = 1250

# Row 27
$ws.Range("D27").Value = 44396
$ws.Range("M27").Value = 12
$ws.Range("N27").Value = 24000
$ws.Range("O27").Value = 24000
$ws.Range("P27").Value = 24000
$ws.Range("S27").Value = 1200

# Row 28
$ws.Range("D28").Value = 44349
$ws.Range("N28").Value = 24000
$ws.Range("O28").Value = 24000
$ws.Range("P28").Value = 24000
$ws.Range("S28").Value = 1200

# Row 29
$ws.Range("D29").Value = 44452
$ws.Range("M29").Value = 25

# Row 30
$ws.Range("D30").Value = 44454
$ws.Range("M30").Value = 25

# Row 31
$ws.Range("D31").Value = 44356
$ws.Range("N31").Value = 24000
$ws.Range("O31").Value = 24000
$ws.Range("P31").Value = 24000
$ws.Range("S31").Value = 1200

# Row 32
$ws.Range("D32").Value = 44469
$ws.Range("M32").Value = 40
$ws.Range("N32").Value = 24000
$ws.Range("O32").Value = 24000
$ws.Range("P32").Value = 24000
$ws.Range("S32").Value = 1200

# Row 33
$ws.Range("D33").Value = 44249
$ws.Range("M33").Value = 15
$ws.Range("N33").Value = 25000
$ws.Range("O33").Value = 25000
$ws.Range("P33").Value = 25000
$ws.Range("S33").Value = 1250

# Row 34
$ws.Range("D34").Value = 44462
$ws.Range("M34").Value = 10

# Row 35
$ws.Range("D35").Value = 44363
$ws.Range("M35").Value = 30
$ws.Range("N35").Value = 24000
$ws.Range("O35").Value = 24000
$ws.Range("P35").Value = 24000
$ws.Range("S35").Value = 1200

# Row 36
$ws.Range("D36").Value = 44474
$ws.Range("M36").Value = 20
$ws.Range("N36").Value = 24000
$ws.Range("O36").Value = 24000
$ws.Range("P36").Value = 24000
$ws.Range("S36").Value = 1200

# Row 37
$ws.Range("D37").Value = 44421
$ws.Range("M37").Value = 20

# Row 38
$ws.Range("D38").Value = 44431
$ws.Range("M38").Value = 40

# Row 39
$ws.Range("D39").Value = 44222
$ws.Range("M39").Value = 15
$ws.Range("N39").Value = 25000
$ws.Range("O39").Value = 25000
$ws.Range("P39").Value = 25000
$ws.Range("S39").Value = 1250

# Row 40
$ws.Range("D40").Value = 44426
$ws.Range("M40").Value = 15
$ws.Range("N40").Value = 24000
$ws.Range("O40").Value = 24000
$ws.Range("P40").Value = 24000
$ws.Range("S40").Value = 1200

# Row 41
$ws.Range("D41").Value = 44425
$ws.Range("M41").Value = 15
$ws.Range("N41").Value = 24000
$ws.Range("O41").Value = 24000
$ws.Range("P41").Value = 24000
$ws.Range("S41").Value = 1200

# Row 42
$ws.Range("D42").Value = 44232
$ws.Range("N42").Value = 25000
$ws.Range("O42").Value = 25000
$ws.Range("P42").Value = 25000
$ws.Range("S42").Value = 1250

# Row 43
$ws.Range("D43").Value = 44468
$ws.Range("M43").Value = 20

# Row 44
$ws.Range("D44").Value = 44238
$ws.Range("M44").Value = 30
$ws.Range("N44").Value = 25000
$ws.Range("O44").Value = 25000
$ws.Range("P44").Value = 25000
$ws.Range("S44").Value = 1250

# Row 45
$ws.Range("D45").Value = 44391
$ws.Range("M45").Value = 10
$ws.Range("N45").Value = 24000
$ws.Range("O45").Value = 24000
$ws.Range("P45").Value = 24000
$ws.Range("S45").Value = 1200

# Row 46
$ws.Range("D46").Value = 44389
$ws.Range("M46").Value = 20

# Row 47
$ws.Range("D47").Value = 44251
$ws.Range("M47").Value = 15
$ws.Range("N47").Value = 25000
$ws.Range("O47").Value = 25000
$ws.Range("P47").Value = 25000
$ws.Range("S47").Value = 1250

# Row 48
$ws.Range("D48").Value = 44432
$ws.Range("M48").Value = 30

# Row 49
$ws.Range("D49").Value = 44400
$ws.Range("M49").Value = 5

